$wb = $excel.ActiveWorkbook

# --- "day" sheet: convert BSE code column D (rows 1175-1188) from text to numeric ---
$dayWs = $wb.Worksheets.Item("day")
$dayWs.Cells.Item(1175, 4).Value = 532977
$dayWs.Cells.Item(1176, 4).Value = 500408
$dayWs.Cells.Item(1177, 4).Value = 500331
$dayWs.Cells.Item(1178, 4).Value = 502355
$dayWs.Cells.Item(1179, 4).Value = 500300
$dayWs.Cells.Item(1180, 4).Value = 524494
$dayWs.Cells.Item(1181, 4).Value = 532689
$dayWs.Cells.Item(1182, 4).Value = 500493
$dayWs.Cells.Item(1183, 4).Value = 500325
$dayWs.Cells.Item(1184, 4).Value = 534816
$dayWs.Cells.Item(1185, 4).Value = 500469
$dayWs.Cells.Item(1186, 4).Value = 541153
$dayWs.Cells.Item(1187, 4).Value = 539437
$dayWs.Cells.Item(1188, 4).Value = 532822

# --- "week" sheet: append 31 new rows (674-704) from the latest stock.yaml "week" pull ---
$weekWs = $wb.Worksheets.Item("week")

$weekWs.Cells.Item(674, 1).Value = 1
$weekWs.Cells.Item(674, 2).Value = 'BOSCHLTD'
$weekWs.Cells.Item(674, 3).Value = 'Bosch Limited'
$weekWs.Cells.Item(674, 4).Formula = "'500530"
$weekWs.Cells.Item(674, 5).Value = -2.69
$weekWs.Cells.Item(674, 6).Value = 34576.95
$weekWs.Cells.Item(674, 7).Value = 19486
$weekWs.Cells.Item(674, 8).Value = "week"
$weekWs.Cells.Item(674, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(675, 1).Value = 2
$weekWs.Cells.Item(675, 2).Value = 'SHREECEM'
$weekWs.Cells.Item(675, 3).Value = 'Shree Cements Limited'
$weekWs.Cells.Item(675, 4).Formula = "'500387"
$weekWs.Cells.Item(675, 5).Value = -2.09
$weekWs.Cells.Item(675, 6).Value = 27041.15
$weekWs.Cells.Item(675, 7).Value = 21938
$weekWs.Cells.Item(675, 8).Value = "week"
$weekWs.Cells.Item(675, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(676, 1).Value = 3
$weekWs.Cells.Item(676, 2).Value = 'ULTRACEMCO'
$weekWs.Cells.Item(676, 3).Value = 'Ultratech Cement Limited'
$weekWs.Cells.Item(676, 4).Formula = "'532538"
$weekWs.Cells.Item(676, 5).Value = -2.13
$weekWs.Cells.Item(676, 6).Value = 11422.8
$weekWs.Cells.Item(676, 7).Value = 299802
$weekWs.Cells.Item(676, 8).Value = "week"
$weekWs.Cells.Item(676, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(677, 1).Value = 4
$weekWs.Cells.Item(677, 2).Value = 'TRENT'
$weekWs.Cells.Item(677, 3).Value = 'Trent Limited'
$weekWs.Cells.Item(677, 4).Formula = "'500251"
$weekWs.Cells.Item(677, 5).Value = -3.67
$weekWs.Cells.Item(677, 6).Value = 6831.55
$weekWs.Cells.Item(677, 7).Value = 768861
$weekWs.Cells.Item(677, 8).Value = "week"
$weekWs.Cells.Item(677, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(678, 1).Value = 5
$weekWs.Cells.Item(678, 2).Value = 'HAL'
$weekWs.Cells.Item(678, 3).Value = 'Hindustan Aeronautics Ltd'
$weekWs.Cells.Item(678, 4).Formula = "'541154"
$weekWs.Cells.Item(678, 5).Value = -4.44
$weekWs.Cells.Item(678, 6).Value = 4190.2
$weekWs.Cells.Item(678, 7).Value = 1598166
$weekWs.Cells.Item(678, 8).Value = "week"
$weekWs.Cells.Item(678, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(679, 1).Value = 6
$weekWs.Cells.Item(679, 2).Value = 'TORNTPHARM'
$weekWs.Cells.Item(679, 3).Value = 'Torrent Pharmaceuticals Limited'
$weekWs.Cells.Item(679, 4).Formula = "'500420"
$weekWs.Cells.Item(679, 5).Value = -0.93
$weekWs.Cells.Item(679, 6).Value = 3437
$weekWs.Cells.Item(679, 7).Value = 358325
$weekWs.Cells.Item(679, 8).Value = "week"
$weekWs.Cells.Item(679, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(680, 1).Value = 7
$weekWs.Cells.Item(680, 2).Value = 'GODREJPROP'
$weekWs.Cells.Item(680, 3).Value = 'Godrej Properties Limited'
$weekWs.Cells.Item(680, 4).Formula = "'533150"
$weekWs.Cells.Item(680, 5).Value = -4.29
$weekWs.Cells.Item(680, 6).Value = 2855.95
$weekWs.Cells.Item(680, 7).Value = 651132
$weekWs.Cells.Item(680, 8).Value = "week"
$weekWs.Cells.Item(680, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(681, 1).Value = 8
$weekWs.Cells.Item(681, 2).Value = 'BALKRISIND'
$weekWs.Cells.Item(681, 3).Value = 'Balkrishna Industries Limited'
$weekWs.Cells.Item(681, 4).Formula = "'502355"
$weekWs.Cells.Item(681, 5).Value = -1.04
$weekWs.Cells.Item(681, 6).Value = 2790.05
$weekWs.Cells.Item(681, 7).Value = 311015
$weekWs.Cells.Item(681, 8).Value = "week"
$weekWs.Cells.Item(681, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(682, 1).Value = 9
$weekWs.Cells.Item(682, 2).Value = 'DEEPAKNTR'
$weekWs.Cells.Item(682, 3).Value = 'Deepak Nitrite Limited'
$weekWs.Cells.Item(682, 4).Formula = "'506401"
$weekWs.Cells.Item(682, 5).Value = -1.59
$weekWs.Cells.Item(682, 6).Value = 2596.85
$weekWs.Cells.Item(682, 7).Value = 225162
$weekWs.Cells.Item(682, 8).Value = "week"
$weekWs.Cells.Item(682, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(683, 1).Value = 10
$weekWs.Cells.Item(683, 2).Value = 'GRASIM'
$weekWs.Cells.Item(683, 3).Value = 'Grasim Industries Limited'
$weekWs.Cells.Item(683, 4).Formula = "'500300"
$weekWs.Cells.Item(683, 5).Value = -1.98
$weekWs.Cells.Item(683, 6).Value = 2488.7
$weekWs.Cells.Item(683, 7).Value = 734974
$weekWs.Cells.Item(683, 8).Value = "week"
$weekWs.Cells.Item(683, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(684, 1).Value = 11
$weekWs.Cells.Item(684, 2).Value = 'TVSMOTOR'
$weekWs.Cells.Item(684, 3).Value = 'Tvs Motor Company Limited'
$weekWs.Cells.Item(684, 4).Formula = "'532343"
$weekWs.Cells.Item(684, 5).Value = -2.76
$weekWs.Cells.Item(684, 6).Value = 2391.65
$weekWs.Cells.Item(684, 7).Value = 575621
$weekWs.Cells.Item(684, 8).Value = "week"
$weekWs.Cells.Item(684, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(685, 1).Value = 12
$weekWs.Cells.Item(685, 2).Value = 'MUTHOOTFIN'
$weekWs.Cells.Item(685, 3).Value = 'Muthoot Finance Limited'
$weekWs.Cells.Item(685, 4).Formula = "'533398"
$weekWs.Cells.Item(685, 5).Value = -2.86
$weekWs.Cells.Item(685, 6).Value = 2075.65
$weekWs.Cells.Item(685, 7).Value = 414899
$weekWs.Cells.Item(685, 8).Value = "week"
$weekWs.Cells.Item(685, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(686, 1).Value = 13
$weekWs.Cells.Item(686, 2).Value = 'DALBHARAT'
$weekWs.Cells.Item(686, 3).Value = 'Dalmia Bharat Limited'
$weekWs.Cells.Item(686, 4).Formula = "'533309"
$weekWs.Cells.Item(686, 5).Value = -2.7
$weekWs.Cells.Item(686, 6).Value = 1801.9
$weekWs.Cells.Item(686, 7).Value = 217114
$weekWs.Cells.Item(686, 8).Value = "week"
$weekWs.Cells.Item(686, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(687, 1).Value = 14
$weekWs.Cells.Item(687, 2).Value = 'PRESTIGE'
$weekWs.Cells.Item(687, 3).Value = 'Prestige Estates Projects Limited'
$weekWs.Cells.Item(687, 4).Formula = "'533274"
$weekWs.Cells.Item(687, 5).Value = -3.72
$weekWs.Cells.Item(687, 6).Value = 1794.5
$weekWs.Cells.Item(687, 7).Value = 2861132
$weekWs.Cells.Item(687, 8).Value = "week"
$weekWs.Cells.Item(687, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(688, 1).Value = 15
$weekWs.Cells.Item(688, 2).Value = 'KOTAKBANK'
$weekWs.Cells.Item(688, 3).Value = 'Kotak Mahindra Bank Limited'
$weekWs.Cells.Item(688, 4).Formula = "'500247"
$weekWs.Cells.Item(688, 5).Value = -1.05
$weekWs.Cells.Item(688, 6).Value = 1743.55
$weekWs.Cells.Item(688, 7).Value = 14686112
$weekWs.Cells.Item(688, 8).Value = "week"
$weekWs.Cells.Item(688, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(689, 1).Value = 16
$weekWs.Cells.Item(689, 2).Value = 'HAVELLS'
$weekWs.Cells.Item(689, 3).Value = 'Havells India Limited'
$weekWs.Cells.Item(689, 4).Formula = "'517354"
$weekWs.Cells.Item(689, 5).Value = -2.06
$weekWs.Cells.Item(689, 6).Value = 1658.25
$weekWs.Cells.Item(689, 7).Value = 735134
$weekWs.Cells.Item(689, 8).Value = "week"
$weekWs.Cells.Item(689, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(690, 1).Value = 17
$weekWs.Cells.Item(690, 2).Value = 'CHOLAFIN'
$weekWs.Cells.Item(690, 3).Value = 'Cholamandalam Investment And Finance Company Limited'
$weekWs.Cells.Item(690, 4).Formula = "'511243"
$weekWs.Cells.Item(690, 5).Value = -1.99
$weekWs.Cells.Item(690, 6).Value = 1189.55
$weekWs.Cells.Item(690, 7).Value = 2300317
$weekWs.Cells.Item(690, 8).Value = "week"
$weekWs.Cells.Item(690, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(691, 1).Value = 18
$weekWs.Cells.Item(691, 2).Value = 'ZYDUSLIFE'
$weekWs.Cells.Item(691, 3).Value = 'Zydus Lifesciences Ltd'
$weekWs.Cells.Item(691, 4).Formula = "'532321"
$weekWs.Cells.Item(691, 5).Value = -1.23
$weekWs.Cells.Item(691, 6).Value = 973.5
$weekWs.Cells.Item(691, 7).Value = 910631
$weekWs.Cells.Item(691, 8).Value = "week"
$weekWs.Cells.Item(691, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(692, 1).Value = 19
$weekWs.Cells.Item(692, 2).Value = 'JINDALSTEL'
$weekWs.Cells.Item(692, 3).Value = 'Jindal Steel & Power Limited'
$weekWs.Cells.Item(692, 4).Formula = "'532286"
$weekWs.Cells.Item(692, 5).Value = -1.56
$weekWs.Cells.Item(692, 6).Value = 908.05
$weekWs.Cells.Item(692, 7).Value = 1099720
$weekWs.Cells.Item(692, 8).Value = "week"
$weekWs.Cells.Item(692, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(693, 1).Value = 20
$weekWs.Cells.Item(693, 2).Value = 'DLF'
$weekWs.Cells.Item(693, 3).Value = 'Dlf Limited'
$weekWs.Cells.Item(693, 4).Formula = "'532868"
$weekWs.Cells.Item(693, 5).Value = -3.9
$weekWs.Cells.Item(693, 6).Value = 830.7
$weekWs.Cells.Item(693, 7).Value = 2873364
$weekWs.Cells.Item(693, 8).Value = "week"
$weekWs.Cells.Item(693, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(694, 1).Value = 21
$weekWs.Cells.Item(694, 2).Value = 'CGPOWER'
$weekWs.Cells.Item(694, 3).Value = 'CG Power and Industrial Solutions Ltd'
$weekWs.Cells.Item(694, 4).Formula = "'500093"
$weekWs.Cells.Item(694, 5).Value = -4.53
$weekWs.Cells.Item(694, 6).Value = 730.05
$weekWs.Cells.Item(694, 7).Value = 2609720
$weekWs.Cells.Item(694, 8).Value = "week"
$weekWs.Cells.Item(694, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(695, 1).Value = 22
$weekWs.Cells.Item(695, 2).Value = 'JSWENERGY'
$weekWs.Cells.Item(695, 3).Value = 'Jsw Energy Limited'
$weekWs.Cells.Item(695, 4).Formula = "'533148"
$weekWs.Cells.Item(695, 5).Value = -2.22
$weekWs.Cells.Item(695, 6).Value = 669.8
$weekWs.Cells.Item(695, 7).Value = 1314435
$weekWs.Cells.Item(695, 8).Value = "week"
$weekWs.Cells.Item(695, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(696, 1).Value = 23
$weekWs.Cells.Item(696, 2).Value = 'GNFC'
$weekWs.Cells.Item(696, 3).Value = 'Gujarat Narmada Valley Fertilizers And Chemicals Limited'
$weekWs.Cells.Item(696, 4).Formula = "'500670"
$weekWs.Cells.Item(696, 5).Value = -2.7
$weekWs.Cells.Item(696, 6).Value = 583.8
$weekWs.Cells.Item(696, 7).Value = 493071
$weekWs.Cells.Item(696, 8).Value = "week"
$weekWs.Cells.Item(696, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(697, 1).Value = 24
$weekWs.Cells.Item(697, 2).Value = 'APOLLOTYRE'
$weekWs.Cells.Item(697, 3).Value = 'Apollo Tyres Limited'
$weekWs.Cells.Item(697, 4).Formula = "'500877"
$weekWs.Cells.Item(697, 5).Value = 0.37
$weekWs.Cells.Item(697, 6).Value = 531.95
$weekWs.Cells.Item(697, 7).Value = 2669279
$weekWs.Cells.Item(697, 8).Value = "week"
$weekWs.Cells.Item(697, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(698, 1).Value = 25
$weekWs.Cells.Item(698, 2).Value = 'VEDL'
$weekWs.Cells.Item(698, 3).Value = 'Vedanta Limited'
$weekWs.Cells.Item(698, 4).Formula = "'500295"
$weekWs.Cells.Item(698, 5).Value = -3.06
$weekWs.Cells.Item(698, 6).Value = 477.25
$weekWs.Cells.Item(698, 7).Value = 13786961
$weekWs.Cells.Item(698, 8).Value = "week"
$weekWs.Cells.Item(698, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(699, 1).Value = 26
$weekWs.Cells.Item(699, 2).Value = 'TATAPOWER'
$weekWs.Cells.Item(699, 3).Value = 'Tata Power Company Limited'
$weekWs.Cells.Item(699, 4).Formula = "'500400"
$weekWs.Cells.Item(699, 5).Value = -2.79
$weekWs.Cells.Item(699, 6).Value = 401.1
$weekWs.Cells.Item(699, 7).Value = 10205219
$weekWs.Cells.Item(699, 8).Value = "week"
$weekWs.Cells.Item(699, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(700, 1).Value = 27
$weekWs.Cells.Item(700, 2).Value = 'CROMPTON'
$weekWs.Cells.Item(700, 3).Value = 'Crompton Greaves Consumer Electricals Limited'
$weekWs.Cells.Item(700, 4).Formula = "'539876"
$weekWs.Cells.Item(700, 5).Value = -1.83
$weekWs.Cells.Item(700, 6).Value = 388.2
$weekWs.Cells.Item(700, 7).Value = 1090955
$weekWs.Cells.Item(700, 8).Value = "week"
$weekWs.Cells.Item(700, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(701, 1).Value = 28
$weekWs.Cells.Item(701, 2).Value = 'INDUSTOWER'
$weekWs.Cells.Item(701, 3).Value = 'Indus Towers Ltd (Bharti Infratel)'
$weekWs.Cells.Item(701, 4).Formula = "'534816"
$weekWs.Cells.Item(701, 5).Value = -2.68
$weekWs.Cells.Item(701, 6).Value = 337.1
$weekWs.Cells.Item(701, 7).Value = 7046696
$weekWs.Cells.Item(701, 8).Value = "week"
$weekWs.Cells.Item(701, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(702, 1).Value = 29
$weekWs.Cells.Item(702, 2).Value = 'JIOFIN'
$weekWs.Cells.Item(702, 3).Value = 'Jio Financial Services Ltd'
$weekWs.Cells.Item(702, 4).Formula = "'20712"
$weekWs.Cells.Item(702, 5).Value = -2.78
$weekWs.Cells.Item(702, 6).Value = 304.3
$weekWs.Cells.Item(702, 7).Value = 19669669
$weekWs.Cells.Item(702, 8).Value = "week"
$weekWs.Cells.Item(702, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(703, 1).Value = 30
$weekWs.Cells.Item(703, 2).Value = 'LTF'
$weekWs.Cells.Item(703, 3).Value = 'L&T Finance Ltd'
$weekWs.Cells.Item(703, 4).Formula = "'533519"
$weekWs.Cells.Item(703, 5).Value = -3.51
$weekWs.Cells.Item(703, 6).Value = 136.48
$weekWs.Cells.Item(703, 7).Value = 8487963
$weekWs.Cells.Item(703, 8).Value = "week"
$weekWs.Cells.Item(703, 9).Value = "20/12/2024 11:34:01"

$weekWs.Cells.Item(704, 1).Value = 31
$weekWs.Cells.Item(704, 2).Value = 'YESBANK'
$weekWs.Cells.Item(704, 3).Value = 'Yes Bank Limited'
$weekWs.Cells.Item(704, 4).Formula = "'532648"
$weekWs.Cells.Item(704, 5).Value = -2.6
$weekWs.Cells.Item(704, 6).Value = 19.83
$weekWs.Cells.Item(704, 7).Value = 76098569
$weekWs.Cells.Item(704, 8).Value = "week"
$weekWs.Cells.Item(704, 9).Value = "20/12/2024 11:34:01"

